$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '삼성전자 삼성 TV용 리모컨(AA59-00598A)'
$ws.Range("B3").Value = 'LG 정품 스마트 TV 매직 리모컨 OLED 올레드 스마트 겸용 새상품'
$ws.Range("B5").Value = '삼성전자 삼성 스마트 리모컨 (BN59-01259A)'
$ws.Range("B7").Value = '삼성 TV리모컨/LG TV 리모컨 (파브 엑스캔버스)'
$ws.Range("B8").Value = 'TV 리모컨 만능 삼성 엘지 리모콘'
$ws.Range("B10").Value = '삼성전자 삼성 스마트 리모컨 (BN59-01276A)'
$ws.Range("B11").Value = '[정품]TM1640A 삼성TV리모컨/BN59-01259A/LED TV/리모콘'
$ws.Range("B12").Value = '더함 3 in 1 통합 리모컨'
$ws.Range("B13").Value = '(5분 설치) 가정용 스마트 통합 만능 리모컨 tv 에어컨 원격제어 음성인식 헤이홈 허브'
$ws.Range("B15").Value = 'LG유플러스 셋톱박스 만능리모컨'
$ws.Range("B16").Value = 'KT 올레 BTV 만능 LG 유플러스 삼성 TV 리모컨 스카이라이프 통합 스마트 헬로 비 티비'
$ws.Range("B17").Value = '올레 스카이라이프 리모컨'
$ws.Range("B19").Value = 'SK브로드밴드 BTV 셋탑박스 리모컨(IPTV 전용 만능리모컨)'
$ws.Range("B20").Value = 'Btv전용 리모컨'
$ws.Range("B21").Value = '[LG정품] 매직리모컨 인공지능 음성인식 20년식 MR20GA AKB75855501'
$ws.Range("B22").Value = '오성통상 IR-1100 하나로 리모컨'
$ws.Range("B23").Value = '와사비망고 딴트공 매직 리모컨'
$ws.Range("B24").Value = 'KT TV 리모컨'
$ws.Range("B25").Value = '헬로TV 리모컨'
$ws.Range("B26").Value = 'LG전자 LG정품 TV 리모컨 (233)'
$ws.Range("B27").Value = 'LG 삼성 TV리모컨 스마트 티비 텔레비젼 엑스캔버스 파브'
$ws.Range("B29").Value = '중소기업TV /삼성 /LG /만능 /통합 /TV리모컨 CB-2000A'
$ws.Range("B30").Value = '넥스디지탈 리모컨'
$ws.Range("B32").Value = 'SK브로드밴드 셋톱박스 리모컨'
$ws.Range("B33").Value = '엘지 TV전용 정품 리모컨 LG 리모콘 (일본어마킹)'
$ws.Range("B35").Value = 'TI320-DU iptv 셋톱박스 lg u+ TV 유플러스리모콘'
$ws.Range("B36").Value = '벽걸이 부착 통합 리모컨 거치대 소품 정리함 케이스 보관함 수납함 주머니 꽂이 수납포켓'
$ws.Range("B38").Value = 'KT올레 셋톱박스 리모컨'
$ws.Range("B40").Value = '삼성전자 삼성 TV용 리모컨(BN59-01312C)'
$ws.Range("B41").Value = 'LG 정품 올레드 매직 리모컨 모음'
$ws.Range("B42").Value = '한국전자 무설정 리모컨 (HK753)'
$ws.Range("B43").Value = '정품 씨제이 리모컨HD 셋톱박스 CJ리모콘 헬로TV리모콘'
$ws.Range("B44").Value = '실리콘 리모콘 커버 야광 리모컨커버 케이스 보호필름'
$ws.Range("B45").Value = 'LG전자 매직 리모컨 (AKB75375503)'
$ws.Range("B46").Value = '스마트리모컨 무선 만능리모컨 원격제어 리모컨 허브 tv 에어컨 통합 헤이홈 iot'
$ws.Range("B47").Value = '올레KT 리모컨 셋톱박스 기가 지니 GiGA Genie 통합리모컨'
$ws.Range("B48").Value = '텐플 무선 만능리모컨 원격제어 리모컨 허브 스마트리모컨 tv 에어컨 iot 구글홈'
$ws.Range("B49").Value = '이노스 매직 리모컨 TV 와 셋톱박스를 하나로 넷플릭스 / 유튜브 핫키 탑재'
$ws.Range("B50").Value = '무설정 엘지TV 만능 리모컨 삼성 티비 리모콘'
$ws.Range("B52").Value = 'PC용 리모콘 곰플레이어 USB 리모컨 무선 마우스기능'
$ws.Range("B53").Value = '삼성전자 삼성 TV용 리모컨(BN59-01175A)'
$ws.Range("B54").Value = 'NEC 프로젝터 리모컨'
$ws.Range("B55").Value = '대우루컴즈 TV리모컨 LUCOMS TV리모컨+건전지무료'
$ws.Range("B56").Value = '부착형 리모컨거치대 리모컨정리함 멀티수납'
$ws.Range("B57").Value = '올레TV 스카이라이프 리모컨'
$ws.Range("B59").Value = '학습형 만능 리모컨/ TV/ 셋탑박스/ 냉난방기 호환'
$ws.Range("B61").Value = 'TV리모컨 삼성 엘지 만능리모컨 통합리모콘 KT 쿡 U플러스 SK BTV 셋톱박스'
$ws.Range("B62").Value = 'KT리모컨 기가지니리모컨 GIGAGenie리모컨 중고제품'
$ws.Range("B63").Value = 'KT리모콘 셋톱박스 / 스카이라이프 리모콘  중고'
$ws.Range("B64").Value = '스마트라TV리모컨 SMATRA TV리모컨+건전지무료'
$ws.Range("B65").Value = '디엘티 모넥스 중소기업TV 리모컨'
$ws.Range("B66").Value = '인켈TV리모컨 INKEL TV리모컨+건전지무료'
$ws.Range("B68").Value = '삼성전자 삼성 TV,비디오용 리모컨(BN59-00377B)'
$ws.Range("B69").Value = 'KT 쿡 올레 TV용 리모컨'
$ws.Range("B72").Value = 'AA59-00598A BN59-01189C 삼성 정품 LCD 3D LED PDP 스마트 TV 리모컨 리모콘'
$ws.Range("B73").Value = '스카이라이프 정품 리모콘 - 벡셀 건전지 무료'
$ws.Range("B74").Value = '통합 만능 무설정 삼성 LG 엘지 티비 TV KT 올레TV 스카이라이프 셋톱박스 리모컨'
$ws.Range("B76").Value = 'LGTV 통합리모컨 유플러스 셋톱박스만능리모컨 LG 엘지 LG전자 리모콘 U+'
$ws.Range("B78").Value = '헤이홈 스마트 리모컨 허브'
$ws.Range("B79").Value = 'PIMPIN 무설정 통합 리모콘  PB-9085'
$ws.Range("B80").Value = '하이얼TV리모컨 HAIER TV리모컨+건전지무료'
$ws.Range("B81").Value = '삼성전자 삼성 TV용 리모컨(00008E)'
$ws.Range("B83").Value = '삼성 TV모니터 리모컨(전기종 호환)'
$ws.Range("B84").Value = 'LG 삼성 TV 리모컨 리모콘(건전지무료)'
$ws.Range("B85").Value = '[무배]  IP TV 리모컨 보호 항균 케이스 / LG U+ SK 브로드밴드 SJD233'
$ws.Range("B86").Value = 'LGTV 삼성TV 스카이라이프 위성방송 TV셋톱박스 케이티 KT 다와 만능리모컨'
$ws.Range("B87").Value = '삼성전자 삼성 TV용 리모컨(BN59-01302A)'
$ws.Range("B88").Value = '삼성전자 삼성 TV용 리모컨(AA83-00654A)'
$ws.Range("B90").Value = 'KT GIGA Genie 리모컨 기가지니 리모콘 중고제품'
$ws.Range("B91").Value = '(정품)TM1680A 삼성TV스마트리모컨/BN59-01243A'
$ws.Range("B92").Value = '대우루컴즈 TV리모컨 / LUCOMS TV리모컨'
$ws.Range("B93").Value = '필립스 TV 리모컨'
$ws.Range("B94").Value = '삼성전자 삼성 TM1240 TV리모컨 AA59-00739A'
$ws.Range("B95").Value = '기가지니 셋톱박스 리모컨'
$ws.Range("B96").Value = 'LG전자 TV 정품 리모컨(AKB74915348)'
$ws.Range("B98").Value = '삼성전자 삼성 TV용 리모컨(AA59-00577A)'
$ws.Range("B99").Value = '리모컨 정리함 삼성 lg 유플러스 sk btv 브로드밴드 kt 기가지니 리모콘 찾기'
$ws.Range("B100").Value = '[당일배송]리모컨 커버/리모콘 실리콘 커버/ 파손방지 리모컨커버/세척가능'
$ws.Range("B101").Value = '스카이 라이프 정품 리모콘                 -  오리지날 정품 쌍방향 리모콘-'
